$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Globo"
$ws.Cells.Item(2, 2).Value = "RJ TV 2"
$ws.Cells.Item(2, 3).Value = "Agricultura"
$ws.Cells.Item(2, 4).Value = "2025-03-31T19:34"
$ws.Cells.Item(2, 5).Value = "Positivo"
$ws.Cells.Item(2, 6).Value = "Com Nota"
$ws.Cells.Item(2, 7).Value = "teste"
